$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the data block (rows 684:685),
# pushing all existing data (old rows 684-734) down by two rows
# (they become rows 686-736).
$ws.Rows("684:685").Insert()

# Populate the newly inserted row 684 (quality "Primera") with the
# latest weekly price observation.
$ws.Cells.Item(684, 1).Value2  = 3
$ws.Cells.Item(684, 2).Value2  = "Femacal de La Calera"
$ws.Cells.Item(684, 3).Value2  = "Coquimbo"
$ws.Cells.Item(684, 4).Value2  = 44783
$ws.Cells.Item(684, 5).Value2  = 5
$ws.Cells.Item(684, 6).Value2  = 100112008
$ws.Cells.Item(684, 7).Value2  = "Coliflor"
$ws.Cells.Item(684, 8).Value2  = "Sin especificar"
$ws.Cells.Item(684, 9).Value2  = "Primera"
$ws.Cells.Item(684, 10).Value2 = 1970
$ws.Cells.Item(684, 11).Value2 = 900
$ws.Cells.Item(684, 12).Value2 = 1000
$ws.Cells.Item(684, 13).Value2 = 957
$ws.Cells.Item(684, 14).Value2 = "`$/unidad"
$ws.Cells.Item(684, 15).Value2 = "Provincia de Quillota"
$ws.Cells.Item(684, 16).Value2 = 957
$ws.Cells.Item(684, 17).Value2 = 1
$ws.Cells.Item(684, 18).Value2 = "Hortaliza"

# Populate the newly inserted row 685 (quality "Segunda") with the
# latest weekly price observation.
$ws.Cells.Item(685, 1).Value2  = 3
$ws.Cells.Item(685, 2).Value2  = "Femacal de La Calera"
$ws.Cells.Item(685, 3).Value2  = "Coquimbo"
$ws.Cells.Item(685, 4).Value2  = 44783
$ws.Cells.Item(685, 5).Value2  = 5
$ws.Cells.Item(685, 6).Value2  = 100112008
$ws.Cells.Item(685, 7).Value2  = "Coliflor"
$ws.Cells.Item(685, 8).Value2  = "Sin especificar"
$ws.Cells.Item(685, 9).Value2  = "Segunda"
$ws.Cells.Item(685, 10).Value2 = 1000
$ws.Cells.Item(685, 11).Value2 = 700
$ws.Cells.Item(685, 12).Value2 = 700
$ws.Cells.Item(685, 13).Value2 = 700
$ws.Cells.Item(685, 14).Value2 = "`$/unidad"
$ws.Cells.Item(685, 15).Value2 = "Provincia de Quillota"
$ws.Cells.Item(685, 16).Value2 = 700
$ws.Cells.Item(685, 17).Value2 = 1
$ws.Cells.Item(685, 18).Value2 = "Hortaliza"

# Make sure column D keeps its date display/number format for the
# newly inserted rows (style index 2 in the original workbook).
$ws.Range("D684:D685").NumberFormat = $ws.Range("D686").NumberFormat

# Refresh the sheet dimension to cover the two additional rows.
$wb.Worksheets.Item(1).UsedRange | Out-Null
